# Updates the "Evaluation form" sheet: refreshed remarks/answers for the
# webgl homework evaluation, a couple of newly-filled-in remark cells,
# row-height bumps to fit the longer wrapped text, and the view scrolled
# down to where the reviewer left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Evaluation form")

# --- Column H remark / indicator text -------------------------------------

$ws.Range("H20").Value = "Floor model, tribune model, road model, wall model, tree model, finish-line model, red car model, starting-light model."
$ws.Range("H21").Value = "Circle geometry, box geometry, cylinder geometry, sphere geometry, torus geometry, custom angle geometry, rectangle model"
$ws.Range("H22").Value = "Models imported: Eltjo (as the race car driver), the car and the lights."
$ws.Range("H23").Value = "The tree has 3 different colors applied, lighter-brown for the tree trunk, darker brown for the tree ground and green for the tree's leafs (see: src/world/environment/tree.js). The floor and road make use of textures (also material properties) which comes down to a total of 5. The final material property is applied for the finish line,  which is a white color for the white blocks."
$ws.Range("H26").Value = "The car is animated"
$ws.Range("H34").Value = "We think that a racing track with a moving car is pretty cool. Look in the race car when it isn't moving to see even more coolness :)"

# --- Row heights, to fit the updated (longer) wrapped text -----------------

$ws.Rows.Item(20).RowHeight = 45
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 30
$ws.Rows.Item(23).RowHeight = 150
$ws.Rows.Item(34).RowHeight = 45.75

# --- Scroll / selection, so the sheet opens where the reviewer left off ----

$ws.Range("H23").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "Evaluation form updated"
